# Add a new "2021" column (R) to the sheet, mirroring column Q.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R4: year header 2021, same formatting as Q4
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R4").Value = 2021

# R5: data value 42.9, same formatting as Q5
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R5").Value = 42.9

# Move the selection the same way the original file shows (R9 now, was Q9)
$ws.Range("R9").Select()
